$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) column cells to Text format so values like "505.50" or
# "58.899.61" are stored verbatim as strings (matching the source data),
# rather than being auto-coerced into numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "58.899.61"
$ws.Range("E2").Value = "  -3.05%  "
$ws.Range("D3").Value = "2.558.59"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "505.50"
$ws.Range("E5").Value = "  -3.43%  "
$ws.Range("D6").Value = "141.93"
$ws.Range("E6").Value = "  -7.79%  "
$ws.Range("D8").Value = "0.552"
$ws.Range("E8").Value = "  -5.77%  "
$ws.Range("D9").Value = "2.561.08"
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("E11").Value = "  -4.05%  "
$ws.Range("D12").Value = "0.329"
$ws.Range("E12").Value = "  -4.93%  "
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").Value = "3.001.85"
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("D15").Value = "58.818.50"
$ws.Range("E15").Value = "  -3.22%  "
$ws.Range("D16").Value = "20.50"
$ws.Range("E16").Value = "  -4.90%  "
$ws.Range("E17").Value = "  -4.71%  "
$ws.Range("D18").Value = "2.558.21"
$ws.Range("E18").Value = "  -1.79%  "
$ws.Range("D19").Value = "4.50"
$ws.Range("E19").Value = "  -5.45%  "
$ws.Range("D20").Value = "330.90"
$ws.Range("E20").Value = "  -6.89%  "
$ws.Range("D21").Value = "10.02"
$ws.Range("E21").Value = "  -4.94%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "5.91"
$ws.Range("E23").Value = "  -4.44%  "
$ws.Range("D24").Value = "59.39"
$ws.Range("E24").Value = "  -2.70%  "
$ws.Range("D25").Value = "0.404"
$ws.Range("E25").Value = "  -4.95%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -5.35%  "
$ws.Range("E28").Value = "  -8.18%  "
$ws.Range("D29").Value = "6.84"
$ws.Range("E29").Value = "  -7.10%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "149.16"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("D32").Value = "18.48"
$ws.Range("E32").Value = "  -4.66%  "
$ws.Range("E33").Value = "  -4.05%  "
$ws.Range("D34").Value = "5.77"
$ws.Range("E34").Value = "  -7.85%  "
$ws.Range("E35").Value = "  -7.41%  "
$ws.Range("D36").Value = "0.874"
$ws.Range("E36").Value = "  -4.80%  "
$ws.Range("D37").Value = "1.10"
$ws.Range("E37").Value = "  -8.16%  "
$ws.Range("D38").Value = "35.74"
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("E39").Value = "  -9.82%  "
$ws.Range("D40").Value = "285.33"
$ws.Range("E40").Value = "  -3.02%  "
$ws.Range("D41").Value = "1.38"
$ws.Range("E41").Value = "  -7.59%  "
$ws.Range("E42").Value = "  -7.86%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "0.0979"
$ws.Range("E44").Value = "  -3.19%  "
$ws.Range("D45").Value = "0.605"
$ws.Range("E45").Value = "  -2.80%  "
$ws.Range("D46").Value = "0.0528"
$ws.Range("E46").Value = "  -5.42%  "
$ws.Range("D47").Value = "10.33"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").Value = "18.55"
$ws.Range("E48").Value = "  -4.96%  "
$ws.Range("E49").Value = "  -5.17%  "
$ws.Range("D50").Value = "4.51"
$ws.Range("E50").Value = "  -8.24%  "
$ws.Range("D51").Value = "1.885.04"
$ws.Range("E51").Value = "  -3.96%  "
